# Gladiator review edit:
# 1. Remove the "Meta description: ..." paragraph that sits right after the
#    H1 title at the top of the document.
# 2. Move that meta-description content to the very end of the document,
#    splitting it into two paragraphs: a bold title line ("Play Gladiator
#    Free: Detailed Game Review") and an italic summary line (the old
#    "Learn about the gameplay..." sentence), replacing the final
#    "Prompt for DALLE: ..." paragraph.

$d = $word.ActiveDocument

# --- Step 1: delete the "Meta description" paragraph near the top ---------
$metaParaIndex = 0
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    if ($d.Paragraphs.Item($i).Range.Text.StartsWith("Meta description:")) {
        $metaParaIndex = $i
        break
    }
}

if ($metaParaIndex -gt 0) {
    $metaPara = $d.Paragraphs.Item($metaParaIndex)
    $nextPara = $d.Paragraphs.Item($metaParaIndex + 1)
    $deleteRange = $d.Range($metaPara.Range.Start, $nextPara.Range.Start)
    $deleteRange.Delete()
}

# --- Step 2: replace the final "Prompt for DALLE" paragraph with the two --
#             new paragraphs (bold title + italic description) ------------
$lastPara = $d.Paragraphs.Item($d.Paragraphs.Count)

if ($lastPara.Range.Text.StartsWith("Prompt for DALLE")) {
    $wNs = 'xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"'
    $titlePara = "<w:p $wNs><w:r/><w:r><w:rPr><w:b/></w:rPr><w:t>Play Gladiator Free: Detailed Game Review</w:t></w:r></w:p>"
    $descPara  = "<w:p $wNs><w:r/><w:r><w:rPr><w:i/></w:rPr><w:t>Learn about the gameplay, features, bonuses, symbols, and RTP of the Gladiator online slot game. Play for free and read our detailed review.</w:t></w:r></w:p>"
    $newXml = $titlePara + $descPara

    $lastPara.Range.InsertXML($newXml) | Out-Null
}

Write-Output "done"
